$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-18 Tuesday" "2024-06-19 Wednesday"

Replace-Text "572×3=" "825×9="
Replace-Text "795×8=" "301×7="
Replace-Text "542×6=" "686×3="
Replace-Text "739×7=" "390×8="
Replace-Text "980×4=" "613×3="

Replace-Text "781×4=" "941×3="
Replace-Text "177×3=" "295×6="
Replace-Text "195×3=" "966×5="
Replace-Text "111×2=" "170×3="
Replace-Text "388×6=" "611×5="

Replace-Text "578×2=" "161×6="
Replace-Text "887×7=" "289×8="
Replace-Text "937×9=" "471×6="
Replace-Text "582×6=" "101×8="
Replace-Text "136×9=" "672×2="

Replace-Text "985×3=" "743×7="
Replace-Text "590×4=" "185×2="
Replace-Text "333×7=" "351×7="
Replace-Text "108×4=" "357×3="
Replace-Text "529×2=" "737×2="

Replace-Text "101×7=" "925×2="
Replace-Text "232×9=" "443×2="
Replace-Text "814×9=" "305×9="
Replace-Text "806×6=" "766×7="
Replace-Text "270×2=" "750×8="
